$wb = $excel.ActiveWorkbook

# "Repayment Schedule" is the sheet carrying the loan's variable-instalment
# schedule; insert a new (blank) column before the existing "Late" column
# (column N) to make room for the new field, shifting Late / Heading /
# Outstanding one column to the right (N->O, O->P, P->Q).
$ws = $wb.Worksheets.Item("Repayment Schedule")
$ws.Columns("N").Insert()

# Make "Repayment Schedule" the active sheet/tab and move its selection,
# so "Transactions" is no longer the active tab.
$ws.Activate()
$ws.Range("S8").Select() | Out-Null
